$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp (row 1) ---
$ws.Cells.Item(1,1).Value = "Datos actualizados a 2 de Julio de 2020 a las 09:10"

# --- Ucrania (row 37): refreshed case counts ---
$ws.Cells.Item(37,2).Value = 45887
$ws.Cells.Item(37,3).Value = 889
$ws.Cells.Item(37,4).Value = 20053
$ws.Cells.Item(37,5).Value = 24649
$ws.Cells.Item(37,7).Value = 12
$ws.Cells.Item(37,8).Value = 1185

# --- Kazajistan (row 39): refreshed case counts ---
$ws.Cells.Item(39,4).Value = 25533
$ws.Cells.Item(39,5).Value = 16853

# --- Senegal / El Salvador swap order (El Salvador now has more cases) ---
# Row 78 becomes El Salvador with freshly updated figures
$ws.Cells.Item(78,1).Value = "El Salvador"
$ws.Cells.Item(78,2).Value = 7000
$ws.Cells.Item(78,3).Value = 264
$ws.Cells.Item(78,4).Value = 4115
$ws.Cells.Item(78,5).Value = 2694
$ws.Cells.Item(78,6).Value = 0
$ws.Cells.Item(78,7).Value = 9
$ws.Cells.Item(78,8).Value = 191

# Row 79 becomes Senegal, carrying the previous row-78 figures
$ws.Cells.Item(79,1).Value = "Senegal"
$ws.Cells.Item(79,2).Value = 6925
$ws.Cells.Item(79,3).Value = 0
$ws.Cells.Item(79,4).Value = 4545
$ws.Cells.Item(79,5).Value = 2264
$ws.Cells.Item(79,6).Value = 0
$ws.Cells.Item(79,7).Value = 0
$ws.Cells.Item(79,8).Value = 116

# --- Letonia (row 132): refreshed case counts ---
$ws.Cells.Item(132,2).Value = 1122
$ws.Cells.Item(132,3).Value = 1
$ws.Cells.Item(132,4).Value = 988
$ws.Cells.Item(132,5).Value = 104

# --- Georgia (row 138): refreshed case counts ---
$ws.Cells.Item(138,2).Value = 939
$ws.Cells.Item(138,3).Value = 8
$ws.Cells.Item(138,4).Value = 817
$ws.Cells.Item(138,5).Value = 107

# --- Taiwan (row 157): refreshed case counts ---
$ws.Cells.Item(157,2).Value = 448
$ws.Cells.Item(157,3).Value = 1
$ws.Cells.Item(157,5).Value = 3

# --- Laos / Santa Lucia swap order (alphabetical-ish reorder, values identical) ---
$ws.Cells.Item(203,1).Value = "Santa Lucia"
$ws.Cells.Item(203,2).Value = 19
$ws.Cells.Item(203,3).Value = 0
$ws.Cells.Item(203,4).Value = 19
$ws.Cells.Item(203,5).Value = 0
$ws.Cells.Item(203,6).Value = 0
$ws.Cells.Item(203,7).Value = 0
$ws.Cells.Item(203,8).Value = 0

$ws.Cells.Item(204,1).Value = "Laos"
$ws.Cells.Item(204,2).Value = 19
$ws.Cells.Item(204,3).Value = 0
$ws.Cells.Item(204,4).Value = 19
$ws.Cells.Item(204,5).Value = 0
$ws.Cells.Item(204,6).Value = 0
$ws.Cells.Item(204,7).Value = 0
$ws.Cells.Item(204,8).Value = 0
